$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 — shifts existing rows 48-104 down to 49-105.
$ws.Rows("48:48").Insert()

# Fill the newly inserted row 48 with the new weekly record.
$ws.Range("A48").Value = 11
$ws.Range("B48").Value = "Vega Monumental Concepción"
$ws.Range("C48").Value = "Bíobío"
$ws.Range("D48").Value = 44638
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100108
$ws.Range("H48").Value = "Tropicales y subtropicales"
$ws.Range("I48").Value = 100108002
$ws.Range("J48").Value = "Mango"
$ws.Range("K48").Value = "Sin especificar"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 220
$ws.Range("N48").Value = 7000
$ws.Range("O48").Value = 7500
$ws.Range("P48").Value = 7227
$ws.Range("Q48").Value = "`$/bandeja 4 kilos"
$ws.Range("R48").Value = "Ecuador"
$ws.Range("S48").Value = 1807
$ws.Range("T48").Value = 4
